# Project "Sample Project" save: update Rules sheet, cell B11 (the R40 rule's
# "To" bound in row 11) from "R40" to "1" — the value must remain TEXT (a
# shared string), not be auto-coerced to the number 1, and the cell's
# existing style (border etc.) must be left untouched.
#
# Plain assignment of a numeric-looking string (e.g. Value = "1") is
# auto-converted to a number by Excel. Forcing text via NumberFormat="@" or
# a leading apostrophe works, but both stamp a *new* cell style (a distinct
# numFmt, or the "quote prefix" flag) onto the xf actually used by the cell —
# which would spuriously add an extra style record. To avoid that, build the
# text "1" as a formula result in a scratch cell (a formula's string result
# carries no quote-prefix baggage), copy just its value/type over to B11 via
# PasteSpecial (which preserves B11's own pre-existing style), then restore
# the scratch cell to its original content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("C11")
$originalScratchValue = $scratch.Value2   # NB: the .Value getter is unreliable here; .Value2 reads correctly

$scratch.Formula = '="1"'
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)   # xlPasteValues: value + data type only, keeps destination formatting
$scratch.Value2 = $originalScratchValue
$excel.CutCopyMode = $false
